$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing data rows 3-11 (values only, styles unchanged) ---
$ws.Range("C3").Value = 8515
$ws.Range("D3").Value = 81
$ws.Range("E3").Value = 108

$ws.Range("C4").Value = 8515
$ws.Range("D4").Value = 75
$ws.Range("E4").Value = 100

$ws.Range("C5").Value = 8515
$ws.Range("D5").Value = 69
$ws.Range("E5").Value = 92

$ws.Range("B6").Value = 9
$ws.Range("C6").Value = 8515
$ws.Range("D6").Value = 63
$ws.Range("E6").Value = 84

$ws.Range("B7").Value = 9
$ws.Range("C7").Value = 8715
$ws.Range("D7").Value = 54
$ws.Range("E7").Value = 72

$ws.Range("B8").Value = 9
$ws.Range("C8").Value = 8915
$ws.Range("D8").Value = 45
$ws.Range("E8").Value = 60

$ws.Range("B9").Value = 9
$ws.Range("C9").Value = 9115
$ws.Range("D9").Value = 36
$ws.Range("E9").Value = 48

$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 9315
$ws.Range("D10").Value = 27
$ws.Range("E10").Value = 36

$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 9515
$ws.Range("D11").Value = 18
$ws.Range("E11").Value = 24

# --- Row 12: previously blank template row, now gets real data ---
# Copy number/border formatting from row 11 so the style matches the
# "filled" look (borderId 6/7 with applyNumberFormat) instead of the
# blank-template look (borderId 6/7 without applyNumberFormat).
$ws.Range("B11:E11").Copy()
$ws.Range("B12:E12").PasteSpecial(-4122)
$ws.Range("B12").Value = 10
$ws.Range("C12").Value = 9515
$ws.Range("D12").Value = 12
$ws.Range("E12").Value = 16

# --- Row 13: previously blank template row, now gets real data ---
# B13 matches the B-column look of row 11; C13/D13/E13 match the
# E-column ("money/right edge") look of row 11.
$ws.Range("B11").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("E11").Copy()
$ws.Range("C13:E13").PasteSpecial(-4122)
$ws.Range("B13").Value = 13
$ws.Range("C13").Value = 9515
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0

$excel.CutCopyMode = $false

# --- Remove the now-unused blank template rows 14-33 ---
$ws.Rows("14:33").Delete()
